$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select column A (the array of numbers) and delete it, shifting cells left
$ws.Range("A1:A3").Select()
$ws.Range("A1:A3").Delete(-4159)
